$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "'001"
$ws.Range("N2").Value = "2018-12-31 00:00:00"

$ws.Range("O2").Value = 293752048.01
$ws.Range("P2").Value = 91929681.16
$ws.Range("Q2").Value = 2651505.23
$ws.Range("R2").Value = 22.6300258769
$ws.Range("S2").Value = 84135473.90000001
$ws.Range("T2").Value = 83.9170774461
$ws.Range("U2").Value = 56637824.44
$ws.Range("V2").Value = 59.2442056396
$ws.Range("W2").Value = 84716870.87
$ws.Range("X2").Value = 43010263.48
$ws.Range("Y2").Value = 153.613804039
$ws.Range("Z2").Value = 352292.68
$ws.Range("AA2").Value = -64.98645256659999
$ws.Range("AB2").Value = 209035177.14
$ws.Range("AC2").Value = 37.2430540793
$ws.Range("AD2").Value = 43.2704180411
$ws.Range("AE2").Value = 60.6826974564
$ws.Range("AF2").Value = 205.0990134597
$ws.Range("AG2").Value = 28.83958476
